# Update column F (dSF) values as part of "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    4  = 2
    5  = -5
    6  = 3
    7  = -1
    10 = 3
    11 = -1
    12 = 4
    13 = 2
    14 = 2
    15 = 2
    16 = 3
    17 = -1
    19 = -5
    20 = -1
    21 = -6
    22 = -1
    23 = -5
    24 = 3
    25 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
